$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for columns B,C,D,E,F,I,J,K,L,M across rows 2-25
$cols = @("B","C","D","E","F","I","J","K","L","M")

$data = @(
    @(1.02,1.05829982881462,1.062805529706472,1.054319396435029,1.070992466017058,1.049024709399867,1.063291768304709,1.065525881468979,1.057062934024869,1.073690826046026),
    @(1.02,1.059565133249015,1.063843583855984,1.055409297003093,1.072223708514826,1.049423142949388,1.064208366506541,1.066378262494624,1.057965346930625,1.074737511571483),
    @(1.02,1.060383472150324,1.064514882631411,1.056114403460226,1.073020403565456,1.049679519781438,1.064800548694466,1.066928812330584,1.058548542978494,1.075414216735087),
    @(1.02,1.06072740925537,1.064797004348008,1.056410800140437,1.073355335811716,1.04978695752508,1.065049283917758,1.067160026081102,1.05879354639089,1.075698568768252),
    @(1.02,1.060785152456983,1.064844368456689,1.056460564728628,1.073411572562233,1.049804976699489,1.065091034910238,1.067198834007942,1.058834673486236,1.075746304870244),
    @(1.02,1.060388068215194,1.064518652716384,1.056118364043415,1.073024878937642,1.04968095671691,1.064803873161908,1.066931902751055,1.058551817400908,1.075418016788457),
    @(1.02,1.058727528300927,1.063156426956003,1.054687761718182,1.071408571589318,1.049159659752623,1.063601728051216,1.065814155076808,1.057368059783578,1.074044677458174),
    @(1.02,1.055798280867717,1.060752944710536,1.052165776217299,1.068560335779566,1.048230032259354,1.061476289313149,1.063836834619051,1.055276505520702,1.071620232053621),
    @(1.02,1.053843154642078,1.05914846583679,1.05048362545457,1.066661323766711,1.047602810429161,1.060054447838836,1.062513346522864,1.053878261542364,1.070000834847552),
    @(1.02,1.052995981265158,1.058453174997814,1.049755016404325,1.065838957660855,1.047329432996013,1.059437592734008,1.061938989129665,1.053271865792037,1.069298857213029),
    @(1.02,1.052681211232781,1.058194829759054,1.049484342428952,1.065533479511559,1.047227619014622,1.05920828458193,1.061725453232273,1.05304647919293,1.06903799380086),
    @(1.02,1.052748734695843,1.058250249489782,1.049542404584044,1.065599006320379,1.047249470668729,1.059257480170846,1.061771266216671,1.053094831941229,1.069093955249595),
    @(1.02,1.052969964174749,1.058431821809929,1.049732643158233,1.065813707057578,1.04732102252304,1.059418641749398,1.0619213421543,1.053253238218802,1.069277296575547),
    @(1.02,1.053106258650885,1.05854368341586,1.049849850572197,1.065945989235688,1.047365072247827,1.059517914707541,1.062013783129237,1.0533508183855,1.069390243557461),
    @(1.02,1.053899366014892,1.059194598433305,1.050531975884755,1.066715899566037,1.047620915863902,1.060095361258713,1.062551437633748,1.053918485914416,1.070047406416303),
    @(1.02,1.054396700764739,1.059602754003413,1.050959792920843,1.067198820401245,1.047780920824021,1.060457258540783,1.062888350493731,1.054274313869271,1.070459420099469),
    @(1.02,1.054686730932832,1.05984077218973,1.051209309631522,1.067480492275732,1.047874076736605,1.060668232590141,1.063084742858169,1.054481771068648,1.070699666769242),
    @(1.02,1.054785614195624,1.059921921471272,1.051294384758562,1.067576533879623,1.047905811306888,1.060740149930533,1.063151686739307,1.054552493197638,1.070781572202055),
    @(1.02,1.054343347351132,1.059558968172984,1.050913894465241,1.067147008379681,1.047763771633393,1.060418442248569,1.062852215691141,1.054236146348346,1.070415222596185),
    @(1.02,1.052904820177165,1.058378355601791,1.04967662365388,1.065750483469658,1.047299959731523,1.059371188686987,1.061877153901822,1.053206595515115,1.069223310391195),
    @(1.02,1.051999826402817,1.057635574887797,1.048898492299277,1.064872346161674,1.047006783407008,1.058711692253058,1.061262970534694,1.052558440669093,1.068473225980628),
    @(1.02,1.052479632527244,1.058029383425707,1.049311015063496,1.065337872428372,1.047162349891784,1.05906140371632,1.061588667884671,1.052902119656789,1.06887092532078),
    @(1.02,1.05436745563753,1.059578753265882,1.05093463406882,1.067170420025739,1.047771521146759,1.060435982014008,1.062868543833719,1.054253392890033,1.070435193773929),
    @(1.02,1.056555953691597,1.061374675693944,1.052817908345904,1.069296697072003,1.048471676212391,1.06202661899612,1.064348942159046,1.055817898385703,1.072247547588061)
)

$rowIndex = 2
foreach ($rowVals in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowIndex).Value = $rowVals[$i]
    }
    $rowIndex++
}
